$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the source diff.
# D-column cells whose new value parses as a plain number need the cell
# NumberFormat forced to Text ("@") first, otherwise Excel auto-converts
# the assigned string into a numeric cell (losing the original text typing).

$ws.Range("D2").Value = '37.817.49'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '2.082.32'
$ws.Range("E3").Value = '  -0.09%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.36'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.625'
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.75'
$ws.Range("E7").Value = '  -0.52%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.393'
$ws.Range("E9").Value = '  +0.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0788'
$ws.Range("E10").Value = '  -0.22%  '
$ws.Range("E11").Value = '  +3.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.01'
$ws.Range("E12").Value = '  +2.13%  '
$ws.Range("D13").Value = '2.389.05'
$ws.Range("E14").Value = '  +0.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.784'
$ws.Range("E15").Value = '  +1.69%  '
$ws.Range("E16").Value = '  +1.99%  '
$ws.Range("D17").Value = '2.073.41'
$ws.Range("E17").Value = '  -0.45%  '
$ws.Range("D18").Value = '37.740.37'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.49'
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("E21").Value = '  +1.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '230.50'
$ws.Range("E22").Value = '  +0.68%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("E24").Value = '  -0.72%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.94'
$ws.Range("E26").Value = '  +10.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '172.04'
$ws.Range("E27").Value = '  +1.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.136'
$ws.Range("E28").Value = '  -1.68%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.52'
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("E31").Value = '  +1.40%  '
$ws.Range("E32").Value = '  +0.95%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0635'
$ws.Range("E33").Value = '  +1.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.68'
$ws.Range("E34").Value = '  -1.53%  '
$ws.Range("E35").Value = '  -1.32%  '
$ws.Range("E36").Value = '  -0.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.41'
$ws.Range("E37").Value = '  -1.26%  '
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.46'
$ws.Range("E39").Value = '  +0.93%  '
$ws.Range("E40").Value = '  +8.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '102.28'
$ws.Range("E41").Value = '  +3.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0978'
$ws.Range("E42").Value = '  -1.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.94'
$ws.Range("E43").Value = '  -0.76%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.84'
$ws.Range("E44").Value = '  +4.57%  '
$ws.Range("D45").Value = '1.451.06'
$ws.Range("E45").Value = '  -0.57%  '
$ws.Range("E46").Value = '  -0.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.10'
$ws.Range("E48").Value = '  -8.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.34'
$ws.Range("E49").Value = '  -0.87%  '
$ws.Range("E50").Value = '  -1.61%  '
$ws.Range("D51").Value = '2.273.71'
$ws.Range("E51").Value = '  -0.16%  '
